$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("Contact.cpp.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.Find.Execute("bool", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
Write-Output ("Found bool: start=" + $rng.Start + " end=" + $rng.End)

$rng.Text = "void"
Write-Output "set text"

$para = $d.Range(2380, 2460).Paragraphs(1)
Write-Output ("para text: [" + $para.Range.Text + "]")

# try adding a bookmark right after "void"
$bmRng = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmRng)
Write-Output "bookmark added"

$checkResult = $para.Range.CheckGrammar()
Write-Output ("CheckGrammar result: " + $checkResult)
$checkResult2 = $para.Range.CheckSpelling()
Write-Output ("CheckSpelling result: " + $checkResult2)
